# Increase size of text: set every paragraph's run formatting (and the
# paragraph mark itself) to 26pt (half-point value 52), matching the
# <w:sz w:val="52"/><w:szCs w:val="52"/> run properties added to both the
# paragraph mark (w:pPr/w:rPr) and each run (w:r/w:rPr) in the target diff.
$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $p.Range.Font.Size = 26
    $p.Range.Font.SizeBi = 26
}
